# Actualización automática 2025-06-16 17:15:08
#
# Applies the June ("junio") sales update for client
# "MOREIRA MOREIRA PATRICIO IGNACIO" (advisor LOZANO MOLINA TITO) across
# the three report sheets: VENTAS POR GRUPO, VENTA MENSUAL and
# CUMPLIMIENTO MENSUAL.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Client's group purchases that were previously unrecorded.
$wsGrupo.Range("C17").Value = 648          # 240X120 PORCELANATO
$wsGrupo.Range("D17").Value = 475.2        # 240X80 PORCELANATO
$wsGrupo.Range("M17").Value = 648.57       # PORCELANATO

# Row 29 tallies "<n> de 27" clients with sales per group; bump the
# counters for the groups that just went from zero to non-zero.
$wsGrupo.Range("C29").Value = "1 de 27"
$wsGrupo.Range("D29").Value = "3 de 27"
$wsGrupo.Range("M29").Value = "1 de 27"

# --- Sheet 2: VENTA MENSUAL -------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F17").Value = 1771.77    # junio - client sale
$wsMensual.Range("F29").Value = 5918.93    # junio - column total

# --- Sheet 3: CUMPLIMIENTO MENSUAL -------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2: 240X120 PORCELANATO
$wsCumpl.Range("D2").Value = 648
$wsCumpl.Range("E2").Value = -303.715395370514
$wsCumpl.Range("F2").Value = 1.882163742690057

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 2344.89
$wsCumpl.Range("E3").Value = 775.2245000000003
$wsCumpl.Range("F3").Value = 0.7515397271478338

# Row 16: PORCELANATO
$wsCumpl.Range("D16").Value = 3327.4
$wsCumpl.Range("E16").Value = 9734.18
$wsCumpl.Range("F16").Value = 0.2547471286015934

# Row 19: TOTAL
$wsCumpl.Range("D19").Value = 12506.69
$wsCumpl.Range("E19").Value = 10993.31093005039
$wsCumpl.Range("F19").Value = 0.5321995534054297
